$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "5    6    7    8    "
$ws.Range("B10").Value = "Q1-"
$ws.Range("B11").Value = "Q1+"
$ws.Range("B12").Value = "Q0-"
$ws.Range("B14").Value = "16  15  14  13"
$ws.Range("B17").Value = "16  15  14  13"
